$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.083.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.27%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.104.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.07%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.100.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.96%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.438"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "

# Row 11
$ws.Range("E11").Value = "  +0.96%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.384"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.34%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.628.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.70%  "

# Row 14
$ws.Range("E14").Value = "  +1.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.46%  "

# Row 16
$ws.Range("E16").Value = "  +1.49%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "58.155.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.095.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.93%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "339.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.87%  "

# Row 23
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.504"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.73%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.15%  "

# Row 26
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0912"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.68%  "

# Row 30
$ws.Range("E30").Value = "  +0.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.10%  "

# Row 32
$ws.Range("E32").Value = "  +2.62%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.85%  "

# Row 34
$ws.Range("E34").Value = "  +2.58%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.62%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.57%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0671"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.140.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.85%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.677"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.53%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "36.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.41%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.296.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.98%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0257"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.33%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.80%  "
